$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (current "FirstName" row). This shifts rows 3-7
# down to 4-8 and - because row 2 above carries the "Hyperlink" cell style in
# column B - Excel's insert logic extends that same formatting onto the new
# B3 automatically, matching the original author's row.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the "ChangedUrl" entry.
$ws.Range("A3").Value = "ChangedUrl"
$ws.Range("B3").Value = "https://www.phptravels.net/account/"

# Move the selection to B3, matching the saved view state.
$ws.Range("B3").Select()
